$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.297.36"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.928.16"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "`'249.04"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "`'0.7160"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "`'0.3208"
$ws.Range("E8").Value = "  -4.15%  "
$ws.Range("D9").Value = "`'27.83"
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("D10").Value = "`'0.07125"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("D11").Value = "`'0.7901"
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").Value = "`'0.07998"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("D13").Value = "1.930.04"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "`'5.386"
$ws.Range("D15").Value = "`'94.85"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "`'14.70"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "30.301.81"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "`'258.12"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").Value = "`'0.000008104"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").Value = "`'5.768"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "2.184.24"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Value = "`'1.001"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").Value = "`'9.539"
$ws.Range("E25").Value = "  -3.34%  "
$ws.Range("D26").Value = "`'164.79"
$ws.Range("E26").Value = "  +2.76%  "
$ws.Range("D27").Value = "`'19.10"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("D28").Value = "`'2.271"
$ws.Range("E28").Value = "  -6.56%  "
$ws.Range("D29").Value = "`'0.1263"
$ws.Range("E29").Value = "  -4.66%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").Value = "`'1.530"
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").Value = "`'4.140"
$ws.Range("E33").Value = "  -2.26%  "
$ws.Range("D34").Value = "`'0.05147"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").Value = "`'1.268"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").Value = "`'0.7439"
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("D39").Value = "`'2.797"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("D40").Value = "`'78.07"
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("D41").Value = "`'6.365"
$ws.Range("E41").Value = "  -4.62%  "
$ws.Range("D42").Value = "`'0.4511"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("D43").Value = "`'1.997"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("D44").Value = "`'0.8488"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "`'9.819"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "`'100.32"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").Value = "`'7.451"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").Value = "`'36.78"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "`'0.4207"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "`'943.63"
$ws.Range("E51").Value = "  +8.77%  "
